$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells with new legend labels
$ws.Range("B1").Value = "large-KDD99"
$ws.Range("C1").Value = "large-CoverType"
$ws.Range("D1").Value = "large-KDD98"

# Update column widths (closest achievable values given engine's internal
# character-width quantization; target stored widths are 13.1640625,
# 16.5 and 12.1640625)
$ws.Columns.Item(2).ColumnWidth = 12.428571428571429
$ws.Columns.Item(3).ColumnWidth = 15.714285714285714
$ws.Columns.Item(4).ColumnWidth = 11.428571428571429

# Update the active selection
$ws.Range("F6").Select()
